$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at the top of this block (rows 195-196), pushing the
# existing rows 195-242 down to 197-244. Excel copies the formatting
# (including the date number format on column D) from the row above on
# insert, matching the style used throughout the sheet.
$ws.Range("A195:A196").EntireRow.Insert()

# New row 195: Sutil De Gase / Primera
$ws.Cells.Item(195, 1).Value = 1
$ws.Cells.Item(195, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(195, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(195, 4).Value = 44711
$ws.Cells.Item(195, 5).Value = 15
$ws.Cells.Item(195, 6).Value = "Fruta"
$ws.Cells.Item(195, 7).Value = 100102
$ws.Cells.Item(195, 8).Value = "Cítricos"
$ws.Cells.Item(195, 9).Value = 100102003
$ws.Cells.Item(195, 10).Value = "Limón"
$ws.Cells.Item(195, 11).Value = "Sutil De Gase"
$ws.Cells.Item(195, 12).Value = "Primera"
$ws.Cells.Item(195, 13).Value = 200
$ws.Cells.Item(195, 14).Value = 34000
$ws.Cells.Item(195, 15).Value = 35000
$ws.Cells.Item(195, 16).Value = 34500
$ws.Cells.Item(195, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(195, 18).Value = "Perú"
$ws.Cells.Item(195, 19).Value = 1438
$ws.Cells.Item(195, 20).Value = 24

# New row 196: Tahití / Primera
$ws.Cells.Item(196, 1).Value = 1
$ws.Cells.Item(196, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(196, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(196, 4).Value = 44711
$ws.Cells.Item(196, 5).Value = 15
$ws.Cells.Item(196, 6).Value = "Fruta"
$ws.Cells.Item(196, 7).Value = 100102
$ws.Cells.Item(196, 8).Value = "Cítricos"
$ws.Cells.Item(196, 9).Value = 100102003
$ws.Cells.Item(196, 10).Value = "Limón"
$ws.Cells.Item(196, 11).Value = "Tahití"
$ws.Cells.Item(196, 12).Value = "Primera"
$ws.Cells.Item(196, 13).Value = 300
$ws.Cells.Item(196, 14).Value = 34000
$ws.Cells.Item(196, 15).Value = 35000
$ws.Cells.Item(196, 16).Value = 34500
$ws.Cells.Item(196, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(196, 18).Value = "Perú"
$ws.Cells.Item(196, 19).Value = 1438
$ws.Cells.Item(196, 20).Value = 24
